$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MS")

# Rename headers
$ws.Range("E1").Value = "Parameter [injection volume setting]"
$ws.Range("G1").Value = "Term Source REF (AFR:0001577)  "
$ws.Range("H1").Value = "Term Accession Number (AFR:0001577)  "

# Clear Term Source REF (MS:1001808) / Term Accession Number (MS:1001808) data cells
$ws.Range("C2:D7").ClearContents()

# Add hyperlinks for Term Accession Number column first (so URL string registered before "UO")
$ws.Hyperlinks.Add($ws.Range("H2"), "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H3"), "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H4"), "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H5"), "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H6"), "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H7"), "http://purl.obolibrary.org/obo/UO_0000101")

# Update Term Source REF value to "UO" for each data row
$ws.Range("G2").Value = "UO"
$ws.Range("G3").Value = "UO"
$ws.Range("G4").Value = "UO"
$ws.Range("G5").Value = "UO"
$ws.Range("G6").Value = "UO"
$ws.Range("G7").Value = "UO"

# Widen column H for the URL text
$ws.Columns("H").ColumnWidth = 43

# Update selection
$ws.Range("E18").Select()
